$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D14").Formula = "=9+16+0"
$ws.Range("E14").Value = 156.77000000000001

$ws.Range("D15").Formula = "=9+16+4"
$ws.Range("E15").Value = 156.77000000000001

$ws.Range("D16").Formula = "=9+16+10"
$ws.Range("E16").Value = 156.77000000000001
